# RP3_FLT_EFF_2022_Jan_Dec.xlsx - "Updates 2022 full year" refresh
#
# 1. Rename the "Belgium" state/FIR label to "Belgium-Luxembourg" in the
#    ERT_FLT_EFF_LOC sheet (A7).
# 2. Bump the "Period Start" / release date on FLT_EFF_YY!B2 from
#    44944 (18 Jan 2023) to 45034 (18 Apr 2023). The other sheets
#    (FLT_EFF_MM, ERT_FLT_EFF_FAB, ERT_FLT_EFF_LOC) hold a formula
#    `=FLT_EFF_YY!B2` in their own B2, so they pick the new date up
#    automatically on recalculation.
# 3. Fill in the previously-blank "D" column (KEA / target) figures for
#    rows 6-33 (Austria..Switzerland) on ERT_FLT_EFF_LOC - except row 7
#    (Belgium-Luxembourg), which stays blank.

$wb = $excel.ActiveWorkbook

$wsYY  = $wb.Worksheets.Item("FLT_EFF_YY")
$wsLOC = $wb.Worksheets.Item("ERT_FLT_EFF_LOC")

# 1) Country label rename
$wsLOC.Range("A7").Value = "Belgium-Luxembourg"

# 2) Refreshed release date (serial date value), FLT_EFF_MM / ERT_FLT_EFF_FAB /
#    ERT_FLT_EFF_LOC reference this cell via formula and recalc automatically.
$wsYY.Range("B2").Value = 45034

# 3) New column D values
$dValues = @{
  6  = 0.0196
  8  = 0.0225
  9  = 0.0146
  10 = 0.0384
  11 = 0.0205
  12 = 0.0114
  13 = 0.0122
  14 = 0.0088
  15 = 0.0283
  16 = 0.023
  17 = 0.0192
  18 = 0.0149
  19 = 0.0113
  20 = 0.0267
  21 = 0.0125
  22 = 0.0192
  23 = 0.018
  24 = 0.0262
  25 = 0.0155
  26 = 0.0165
  27 = 0.018
  28 = 0.0205
  29 = 0.0213
  30 = 0.0155
  31 = 0.0308
  32 = 0.0105
  33 = 0.0395
}

foreach ($row in $dValues.Keys) {
  $wsLOC.Cells.Item($row, 4).Value = $dValues[$row]
}
